# Apply the "Office Theme" design (built-in default theme colours) to the
# presentation's slide master theme, replacing the "Integral" / "Red Violet"
# colour scheme currently applied, and switch the table on slide 5 to the
# corresponding built-in table style.

$p = $ppt.ActivePresentation

# --- 1. Re-colour the deck's theme (Design > Office Theme) ----------------
# ThemeColorScheme items are ordered exactly like <a:clrScheme>:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# RGB() is not available in this host, so colours are packed manually as
# R + G*256 + B*65536 (matching the VBA RGB() long value layout).
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0            # dk1      000000
$cs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388      # dk2      44546A
$cs.Item(4).RGB  = 15132391     # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939     # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501      # accent2  ED7D31
$cs.Item(7).RGB  = 10855845     # accent3  A5A5A5
$cs.Item(8).RGB  = 49407        # accent4  FFC000
$cs.Item(9).RGB  = 12874308     # accent5  4472C4
$cs.Item(10).RGB = 4697456      # accent6  70AD47
$cs.Item(11).RGB = 12673797     # hlink    0563C1
$cs.Item(12).RGB = 7491477      # folHlink 954F72

# --- 2. Re-style the table on slide 5 --------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{8C694D1B-8A29-4697-A5E1-6C1AB84A355D}")
